$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.611
$ws.Range("D4").Value = -7.393000000000001
$ws.Range("A9").Value = -20.912
$ws.Range("D10").Value = -7.907000000000001
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("C21").Value = -12.688
